# Update crypto price (D) and 1h volume change (E) columns with freshly
# scraped values from the GitHub Actions cron run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.119.64"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").Value = "'1.860.72"
$ws.Range("E3").Value = "  -4.52%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'233.39"
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4663"
$ws.Range("E7").Value = "  -3.49%  "
$ws.Range("D8").Value = "'0.2811"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").Value = "'0.06544"
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("D10").Value = "'19.62"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").Value = "'0.07806"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'96.47"
$ws.Range("E12").Value = "  -8.27%  "
$ws.Range("D13").Value = "'1.871.72"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "'5.132"
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").Value = "'0.6663"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "'280.10"
$ws.Range("E16").Value = "  -6.30%  "
$ws.Range("D17").Value = "'30.154.74"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "'0.9994"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'5.498"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").Value = "'2.101.84"
$ws.Range("E21").Value = "  -5.03%  "
$ws.Range("D22").Value = "'0.000007227"
$ws.Range("E22").Value = "  -5.47%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("D24").Value = "'6.124"
$ws.Range("E24").Value = "  -5.65%  "
$ws.Range("D25").Value = "'9.308"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "'165.55"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'18.85"
$ws.Range("E27").Value = "  -5.37%  "
$ws.Range("D28").Value = "'1.908"
$ws.Range("E28").Value = "  -11.32%  "
$ws.Range("D29").Value = "'1.340"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").Value = "'0.09542"
$ws.Range("E30").Value = "  -6.47%  "
$ws.Range("D31").Value = "'4.402"
$ws.Range("E31").Value = "  -5.43%  "
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("D33").Value = "'4.103"
$ws.Range("E33").Value = "  -6.38%  "
$ws.Range("D34").Value = "'0.04651"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("D35").Value = "'0.7010"
$ws.Range("E35").Value = "  -6.32%  "
$ws.Range("D36").Value = "'1.091"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("D37").Value = "'2.708"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "'0.01852"
$ws.Range("E38").Value = "  -6.07%  "
$ws.Range("D39").Value = "'6.278"
$ws.Range("E39").Value = "  -5.85%  "
$ws.Range("D40").Value = "'2.510"
$ws.Range("E40").Value = "  -5.41%  "
$ws.Range("D41").Value = "'73.41"
$ws.Range("E41").Value = "  -5.48%  "
$ws.Range("D42").Value = "'0.8520"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").Value = "'1.914"
$ws.Range("E43").Value = "  -6.35%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'103.66"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").Value = "'0.4147"
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("D47").Value = "'994.48"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("D48").Value = "'7.175"
$ws.Range("E48").Value = "  -6.06%  "
$ws.Range("D49").Value = "'9.307"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "'34.14"
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("D51").Value = "'0.1139"
$ws.Range("E51").Value = "  -6.65%  "
